# Edit script for LOQ4209.xlsx rebuild (commit: "Build site at 2022-09-26 16:07:08 UTC")
# Reshuffles the course-info table: drops the "Bibliografia" reference text and the
# trailing "Requisitos" value row, shifting remaining labelled rows up so the grid now
# spans A1:C23 instead of A1:C24.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (old "LOB1012 - Estatistica..." requisito value, row 24).
$ws.Rows(24).Delete()

# Rewrite the affected cells (rows 10, 13-23) with the final values.
$ws.Range("B10").Value = '5840535 - Messias Borges Silva'
$ws.Range("C10").Value = '5840535 - Messias Borges Silva'
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = '1-Philosophy, basic concepts. 2 – Total Quality Management Tools. 3 – Quality Management Systems'
$ws.Range("C14").Value = '1-Philosophy, basic concepts. 2 – Total Quality Management Tools. 3 – Quality Management Systems'
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2022'
$ws.Range("C15").Value = '01/01/2022'
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1-PHILOSOPHY, BASIC CONCEPTSDefinition of Quality. Competitiveness. Quality History. Internal and External Benefits of Quality2 – TOTAL QUALITY MANAGEMENT TOOLSRoutine management: MASP Solving Problem Method, Brainstorming, Flowcharts, PDCA and SDCA, Basic Management Units, 5S Program, Operating Procedures and Work Instructions, Internal Audits, Education Program, Visual Management, Customer Monitoring, Kaizen Groups.3 – QUALITY MANAGEMENT SYSTEMSCertification Systems: ISO (9001, 14001, 17025 and 65), SA 8000, OHSAS 18000Implementation methodology, documentation, requirements, participation of staff, middle management and operational instances.Pre-audit, certification audit, maintenance audits.'
$ws.Range("C16").Value = '1-PHILOSOPHY, BASIC CONCEPTSDefinition of Quality. Competitiveness. Quality History. Internal and External Benefits of Quality2 – TOTAL QUALITY MANAGEMENT TOOLSRoutine management: MASP Solving Problem Method, Brainstorming, Flowcharts, PDCA and SDCA, Basic Management Units, 5S Program, Operating Procedures and Work Instructions, Internal Audits, Education Program, Visual Management, Customer Monitoring, Kaizen Groups.3 – QUALITY MANAGEMENT SYSTEMSCertification Systems: ISO (9001, 14001, 17025 and 65), SA 8000, OHSAS 18000Implementation methodology, documentation, requirements, participation of staff, middle management and operational instances.Pre-audit, certification audit, maintenance audits.'
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5840535 - Messias Borges Silva'
$ws.Range("C18").Value = '5840535 - Messias Borges Silva'
$ws.Range("A19").Value = 'Critério:'
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = 'LOB1012 -  Estatística  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOB1012 -  Estatística  (Requisito fraco)
'

# Adjust row heights to match the new layout.
$ws.Rows(13).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(17).AutoFit()
$ws.Rows(18).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(22).AutoFit()
$ws.Rows(23).RowHeight = 30
